$wb = $excel.ActiveWorkbook
$wbs = $wb.Worksheets.Item("WBS")
$wbs.Activate()
$win = $excel.ActiveWindow
$sv = $win.SheetViews
Write-Output $sv
$members = $sv | Get-Member
Write-Output $members
